$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N5").Value = "Daniel, Emily"
$ws.Range("B2:B5").ClearFormats()
$ws.Range("O10").Select()
